# Fruta / hortaliza, semanal
# Insert a new weekly record as row 44 (Angeleno, Primera, Provincia de Curicó),
# pushing the existing rows 44-61 down to 45-62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44; everything below shifts down by one.
$ws.Rows.Item(44).Insert()

# Populate the new row 44 with the new record's data.
$ws.Cells.Item(44, 1).Value = 11
$ws.Cells.Item(44, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(44, 3).Value = "Bíobío"
$ws.Cells.Item(44, 4).Value = 44637
$ws.Cells.Item(44, 5).Value = 8
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100103
$ws.Cells.Item(44, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(44, 9).Value = 100103002
$ws.Cells.Item(44, 10).Value = "Ciruela"
$ws.Cells.Item(44, 11).Value = "Angeleno"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 220
$ws.Cells.Item(44, 14).Value = 8000
$ws.Cells.Item(44, 15).Value = 8500
$ws.Cells.Item(44, 16).Value = 8273
$ws.Cells.Item(44, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(44, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(44, 19).Value = 460
$ws.Cells.Item(44, 20).Value = 18
